# Update NIFTY50 sheet data: overwrite rows 2-23 with refreshed data
# and append two new trailing rows (24-25). Mirrors the author's
# "updated code with new" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A=45352; B="22116.30"; C="22061.60"; D="22131.50"; E="22061.60"; F="0"; G="0.37%" },
    @{ Row=3; A=45323; B="22034.10"; C="21774.05"; D="22293.20"; E="21531.90"; F="6368.00B"; G="1.42%" },
    @{ Row=4; A=45292; B="21725.70"; C="21727.75"; D="22124.15"; E="21137.20"; F="6824.38B"; G="-0.03%" },
    @{ Row=5; A=45261; B="21731.40"; C="20194.10"; D="21801.45"; E="20183.70"; F="6121.37B"; G="7.94%" },
    @{ Row=6; A=45231; B="20133.15"; C="19064.05"; D="20158.70"; E="18973.70"; F="4423.32B"; G="5.52%" },
    @{ Row=7; A=45200; B="19079.60"; C="19622.40"; D="19849.75"; E="18837.85"; F="4248.51B"; G="-2.84%" },
    @{ Row=8; A=45170; B="19638.30"; C="19258.15"; D="20222.45"; E="19255.70"; F="5666.41B"; G="2.00%" },
    @{ Row=9; A=45139; B="19253.80"; C="19784.00"; D="19795.60"; E="19223.65"; F="6253.99B"; G="-2.53%" },
    @{ Row=10; A=45108; B="19753.80"; C="19246.50"; D="19991.85"; E="19234.40"; F="5802.27B"; G="2.94%" },
    @{ Row=11; A=45078; B="19189.05"; C="18579.40"; D="19201.70"; E="18464.55"; F="5144.19B"; G="3.53%" },
    @{ Row=12; A=45047; B="18534.40"; C="18124.80"; D="18662.45"; E="18042.40"; F="5737.40B"; G="2.60%" },
    @{ Row=13; A=45017; B="18065.00"; C="17427.95"; D="18089.15"; E="17312.75"; F="4459.76B"; G="4.06%" },
    @{ Row=14; A=44986; B="17359.75"; C="17360.10"; D="17799.95"; E="16828.35"; F="5622.28B"; G="0.32%" },
    @{ Row=15; A=44958; B="17303.95"; C="17811.60"; D="18134.75"; E="17255.20"; F="5685.63B"; G="-2.03%" },
    @{ Row=16; A=44927; B="17662.15"; C="18131.70"; D="18251.95"; E="17405.55"; F="5632.81B"; G="-2.45%" },
    @{ Row=17; A=44896; B="18105.30"; C="18871.95"; D="18887.60"; E="17774.25"; F="4741.44B"; G="-3.48%" },
    @{ Row=18; A=44866; B="18758.35"; C="18130.70"; D="18816.05"; E="17959.20"; F="5257.00B"; G="4.14%" },
    @{ Row=19; A=44835; B="18012.20"; C="17102.10"; D="18022.80"; E="16855.55"; F="4539.96B"; G="5.37%" },
    @{ Row=20; A=44805; B="17094.35"; C="17485.70"; D="18096.15"; E="16747.70"; F="6896.51B"; G="-3.74%" },
    @{ Row=21; A=44774; B="17759.30"; C="17243.20"; D="17992.20"; E="17154.80"; F="5589.51B"; G="3.50%" },
    @{ Row=22; A=44743; B="17158.25"; C="15703.70"; D="17172.80"; E="15511.05"; F="5475.27B"; G="8.73%" },
    @{ Row=23; A=44713; B="15780.25"; C="16594.40"; D="16793.85"; E="15183.40"; F="5514.20B"; G="-4.85%" },
    @{ Row=24; A=44682; B="16584.55"; C="16924.45"; D="17132.85"; E="15735.75"; F="6343.31B"; G="-3.03%" },
    @{ Row=25; A=44652; B="17102.55"; C="17436.90"; D="18114.65"; E="16824.70"; F="5658.25B"; G="-2.07%" }
)

foreach ($item in $rows) {
    $aAddr = "A" + $item.Row
    $ws.Range($aAddr).NumberFormat = "YYYY-MM-DD"
    $ws.Range($aAddr).Value = $item.A

    foreach ($col in @("B", "C", "D", "E", "F", "G")) {
        $addr = $col + $item.Row
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $item[$col]
        $ws.Range($addr).Style = "Normal"
    }
}
